# Implement PDF download function
# Adds a new "[bill_company_name]" placeholder to the bill template sheet,
# right after the existing "[bill_payment_type]" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new placeholder value to the next free row (row 12 -> row 13).
$ws.Range("A13").Value = "[bill_company_name]"

# Leave the active selection on the following empty row, as in the source edit.
$ws.Range("A14").Select() | Out-Null
